$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update header row (A1:J1) ---
# Columns: model, level, auc_roc, J, Sens, Spec, Acc, BA, ap, MCC
# Copy the existing bold/bordered header style (from D1) onto the new header cells E1:J1
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$headers = New-Object 'object[,]' 1,10
$headers[0,0] = "model"
$headers[0,1] = "level"
$headers[0,2] = "auc_roc"
$headers[0,3] = "J"
$headers[0,4] = "Sens"
$headers[0,5] = "Spec"
$headers[0,6] = "Acc"
$headers[0,7] = "BA"
$headers[0,8] = "ap"
$headers[0,9] = "MCC"
$ws.Range("A1:J1").Value2 = $headers

# --- 2) Replace the data rows (A2:J12) with the recalculated metrics ---
$data = New-Object 'object[,]' 11,10
$data[0,0] = "TAP"
$data[0,1] = 10
$data[0,2] = 0.8281573498964804
$data[0,3] = 0.2604
$data[0,4] = 0.5908289241622575
$data[0,5] = 0.9347826086956522
$data[0,6] = 0.6166394779771615
$data[0,7] = 0.7628057664289549
$data[0,8] = 0.9829796549385335
$data[0,9] = 0.2784265423272125
$data[1,0] = "TAP"
$data[1,1] = 12
$data[1,2] = 0.8288523582731602
$data[1,3] = 0.1282
$data[1,4] = 0.6554770318021201
$data[1,5] = 0.8913043478260869
$data[1,6] = 0.673202614379085
$data[1,7] = 0.7733906898141035
$data[1,8] = 0.9827421423411445
$data[1,9] = 0.2961767232155829
$data[2,0] = "TAP"
$data[2,1] = 11
$data[2,2] = 0.8246325690770135
$data[2,3] = 0.1481
$data[2,4] = 0.6384479717813051
$data[2,5] = 0.8888888888888888
$data[2,6] = 0.6568627450980392
$data[2,7] = 0.7636684303350969
$data[2,8] = 0.9828679887234953
$data[2,9] = 0.2809117087133983
$data[3,0] = "TAP"
$data[3,1] = 6
$data[3,2] = 0.8231347289318304
$data[3,3] = 0.0234
$data[3,4] = 0.7160493827160493
$data[3,5] = 0.8260869565217391
$data[3,6] = 0.7243066884176182
$data[3,7] = 0.7710681696188942
$data[3,8] = 0.9828354738421248
$data[3,9] = 0.3050368134163043
$data[4,0] = "TAP"
$data[4,1] = 7
$data[4,2] = 0.8279656468062264
$data[4,3] = 0.1416
$data[4,4] = 0.6455026455026455
$data[4,5] = 0.8913043478260869
$data[4,6] = 0.6639477977161501
$data[4,7] = 0.7684034966643662
$data[4,8] = 0.9832150769236028
$data[4,9] = 0.2893303367574673
$data[5,0] = "TAP"
$data[5,1] = 3
$data[5,2] = 0.8162717583007437
$data[5,3] = 0.162
$data[5,4] = 0.6843033509700176
$data[5,5] = 0.8913043478260869
$data[5,6] = 0.6998368678629691
$data[5,7] = 0.7878038493980523
$data[5,8] = 0.9826767333598202
$data[5,9] = 0.3161472611002694
$data[6,0] = "TAP"
$data[6,1] = 9
$data[6,2] = 0.8267770876466529
$data[6,3] = 0.1633
$data[6,4] = 0.6296296296296297
$data[6,5] = 0.8913043478260869
$data[6,6] = 0.6492659053833605
$data[6,7] = 0.7604669887278583
$data[6,8] = 0.9829380728290058
$data[6,9] = 0.2791013801743985
$data[7,0] = "TAP"
$data[7,1] = 8
$data[7,2] = 0.82587608312246
$data[7,3] = 0.176
$data[7,4] = 0.6278659611992945
$data[7,5] = 0.9130434782608695
$data[7,6] = 0.6492659053833605
$data[7,7] = 0.770454719730082
$data[7,8] = 0.9828899491286367
$data[7,9] = 0.2894566254751216
$data[8,0] = "TAP"
$data[8,1] = 4
$data[8,2] = 0.8156966490299824
$data[8,3] = 0.1124
$data[8,4] = 0.6649029982363316
$data[8,5] = 0.8695652173913043
$data[8,6] = 0.6802610114192496
$data[8,7] = 0.767234107813818
$data[8,8] = 0.9824309618771062
$data[8,9] = 0.2908231830571419
$data[9,0] = "TAP"
$data[9,1] = 2
$data[9,2] = 0.7945326278659611
$data[9,3] = 0.2051
$data[9,4] = 0.6261022927689595
$data[9,5] = 0.8043478260869565
$data[9,6] = 0.6394779771615008
$data[9,7] = 0.715225059427958
$data[9,8] = 0.9798420105745118
$data[9,9] = 0.2309100907632188
$data[10,0] = "TAP"
$data[10,1] = 5
$data[10,2] = 0.8165018020090485
$data[10,3] = 0.1761
$data[10,4] = 0.6278659611992945
$data[10,5] = 0.9347826086956522
$data[10,6] = 0.6508972267536705
$data[10,7] = 0.7813242849474733
$data[10,8] = 0.9823820912481167
$data[10,9] = 0.300914823355256

$ws.Range("A2:J12").Value2 = $data

Write-Host "Edit complete"
